# NSantiam HBVCALIB and gages - new spinup after extending DET12 down to the
# Niagara gage location, and adjustment of the gage area:HBVCALIB area factors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the Breitenbush/Boulder Creek gage rows:
#   row 14 will hold a new "Nsantiam" totals line,
#   row 15 stays as a blank spacer row (matches the existing blank row 14 below row 13).
$ws.Rows("14:15").Insert() | Out-Null

# DET12's pour-point drainage area grew now that it extends down to the Niagara gage.
$ws.Range("D11").Value = 1068300000

# New totals row: name + summed area (m2) + area (km2).
$ws.Range("B14").Value = "Nsantiam"
$ws.Range("D14").Formula = "=SUM(D9:D13)"

# Extend the area-in-km2 shared formula down through the new total row and
# through the (re-numbered) Breitenbush/Boulder Creek rows.
$ws.Range("E10:E14").Formula = "=D10/1000000"
$ws.Range("E17:E18").Formula = "=D17/1000000"

# Likewise for the gage-area-in-km2 formula on the (re-numbered) gage rows.
$ws.Range("P10:P13").Formula = "=O10*2.58998811"
$ws.Range("P17:P18").Formula = "=O17*2.58998811"

# The spacer row (15) shouldn't carry the area formula/value that Excel
# copied down automatically when the rows were inserted.
$ws.Range("D15").Clear()

# Leave the selection where the author last left it.
$ws.Range("E14").Select() | Out-Null
